$d = $word.ActiveDocument

# Remove the existing "_GoBack" bookmark (currently wraps the title picture paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-add the "_GoBack" bookmark around the empty paragraph in the
# "Use Case Overview:" value cell (Table 1, row 6, column 2).
$cellRange = $d.Tables(1).Cell(6, 2).Range
$d.Bookmarks.Add("_GoBack", $cellRange)
